$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update column F ("想去人数") values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 0
$wsExhibit.Range("F5").Value = 201
$wsExhibit.Range("F6").Value = 0
$wsExhibit.Range("F7").Value = 0
$wsExhibit.Range("F10").Value = 730
$wsExhibit.Range("F11").Value = 213
$wsExhibit.Range("F12").Value = 1113
$wsExhibit.Range("F13").Value = 98
$wsExhibit.Range("F15").Value = 0
$wsExhibit.Range("F16").Value = 78
$wsExhibit.Range("F17").Value = 137
$wsExhibit.Range("F19").Value = 0
$wsExhibit.Range("F20").Value = 6153
$wsExhibit.Range("F22").Value = 0
$wsExhibit.Range("F25").Value = 46
$wsExhibit.Range("F26").Value = 0
$wsExhibit.Range("F27").Value = 0
$wsExhibit.Range("F28").Value = 34
$wsExhibit.Range("F29").Value = 2532
$wsExhibit.Range("F31").Value = 0
$wsExhibit.Range("F33").Value = 262
$wsExhibit.Range("F36").Value = 157
$wsExhibit.Range("F37").Value = 1552
$wsExhibit.Range("F39").Value = 0
$wsExhibit.Range("F40").Value = 51
$wsExhibit.Range("F42").Value = 0
$wsExhibit.Range("F43").Value = 0
$wsExhibit.Range("F44").Value = 0
$wsExhibit.Range("F45").Value = 568

# Sheet "演出" (sheet2) - update column F value
$wsPerform = $wb.Worksheets.Item("演出")
$wsPerform.Range("F2").Value = 106

# Sheet "全部类型" (sheet4) - update column F ("想去人数") values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 0
$wsAll.Range("F3").Value = 218
$wsAll.Range("F7").Value = 0
$wsAll.Range("F8").Value = 0
$wsAll.Range("F9").Value = 0
$wsAll.Range("F10").Value = 91
$wsAll.Range("F11").Value = 0
$wsAll.Range("F12").Value = 213
$wsAll.Range("F13").Value = 1113
$wsAll.Range("F14").Value = 98
$wsAll.Range("F16").Value = 170
$wsAll.Range("F18").Value = 0
$wsAll.Range("F20").Value = 0
$wsAll.Range("F21").Value = 6153
$wsAll.Range("F22").Value = 39
$wsAll.Range("F23").Value = 38
$wsAll.Range("F24").Value = 0
$wsAll.Range("F25").Value = 531
$wsAll.Range("F27").Value = 3940
$wsAll.Range("F28").Value = 391
$wsAll.Range("F30").Value = 2532
$wsAll.Range("F35").Value = 0
$wsAll.Range("F37").Value = 157
$wsAll.Range("F38").Value = 1552
$wsAll.Range("F39").Value = 933
$wsAll.Range("F40").Value = 40
$wsAll.Range("F41").Value = 0
$wsAll.Range("F43").Value = 481
$wsAll.Range("F44").Value = 0
$wsAll.Range("F46").Value = 568
